# regen sval data to filter save games
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(3.182878228561681, 1.65323645889881, 0.1529057820181812, 0.4998867070740569, 5.488907176552729),
    @(1.505614041169197, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 4.371470058157054),
    @(3.182878228561681, 1.65323645889881, 3.082599426703578, 6.48142807727062, 14.40014219143469),
    @(0.1554434735375247, 0.3375848360084654, 3.082599426703578, 6.48142807727062, 10.05705581352019),
    @(3.182878228561681, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 6.048734245549538),
    @(0.1554434735375247, 0.3375848360084654, 0.7127328510149897, 6.48142807727062, 7.6871892378316),
    @(3.182878228561681, 1.65323645889881, 0.1529057820181812, 0.4998867070740569, 5.488907176552729),
    @(3.182878228561681, 1.65323645889881, 16.98373111632243, 6.48142807727062, 28.30127388105354),
    @(0.3464964993005633, 0.05231270169004087, 0.1529057820181812, 0.4998867070740569, 1.051601690082842)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $vals = $data[$i]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 7).Value = $vals[4]
}
